$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.919.10"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "1.890.14"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +1.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.53"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.015"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4693"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3929"
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.82"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08061"
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.022"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.84"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "1.880.76"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.969"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.131"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06791"
$ws.Range("E17").Value = "  +3.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001050"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "87.32"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.22"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").Value = "27.919.35"
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.509"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.01"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.336"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").Value = "2.103.25"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.51"
$ws.Range("E27").Value = "  +3.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.08"
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.083"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.478"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.13"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9717"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09505"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.644"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.407"
$ws.Range("E35").Value = "  -4.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.377"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06125"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02256"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.218"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.048"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5985"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1884"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.30"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.265"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5708"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.23"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.404"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.935"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06934"
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.00"
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.071"
$ws.Range("E51").Value = "  +0.57%  "
